$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 40, shifting existing rows 40:154 down to 41:155
$ws.Rows("40:40").Insert()

# Populate the newly inserted row 40 with the new entry
$ws.Range("R40").Value = "communication feedback"
$ws.Range("S40").Value = "2024-09-16 11:13:15"
